# edit.ps1 - Applies the Notes.docx revision described in the commit:
# "add addtional condition for comparing attendee company."
#
# The script performs the textual / structural edits captured by the
# target XML diff using the Word COM object model surface exposed by
# $word.ActiveDocument ($d below).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...OCCCIO 2019." -> "...OCCCIO 2020."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "It has passed down to Durham College for further development for OCCCIO 2019.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "It has passed down to Durham College for further development for OCCCIO 2020.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) "...generate a password for the Admin user..." ->
#    "...generate a password for the Admin or staff user..."
#    plus a _GoBack bookmark right after the inserted text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Next, you need to generate a password for the Admin user, this is done using the",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Next, you need to generate a password for the Admin or staff user, this is done using the",
    2) | Out-Null

# Place the _GoBack bookmark right after " or staff" (collapsed range).
$bookmarkRange = $d.Content
$bookmarkRange.Find.Execute("Admin or staff", $true) | Out-Null
$bookmarkRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ---------------------------------------------------------------------
# 3) Insert two new paragraphs after the "Alter the Applications
#    Arguments..." paragraph.
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute('Alter the Applications Arguments to "p [path]".', $true) | Out-Null
$anchorPara = $anchor.Paragraphs(1)
$insertPos = $anchorPara.Range.End
$anchorPara.Range.InsertParagraphAfter() | Out-Null

$examplePara = $d.Range($insertPos, $insertPos).Paragraphs(1)
$examplePara.Range.InsertAfter("Example: ") | Out-Null
$examplePara.Range.InsertAfter("p C:/PrizeDraw") | Out-Null

$insertPos2 = $examplePara.Range.End
$examplePara.Range.InsertParagraphAfter() | Out-Null
$genPara = $d.Range($insertPos2, $insertPos2).Paragraphs(1)
$genPara.Range.InsertAfter("Will generate password on prizedraw folder located on c drive") | Out-Null

# ---------------------------------------------------------------------
# 4) Typo fix: "ahve" -> "have"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("ahve", $true, $false, $false, $false, $false,
                         $true, 1, $false, "have", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Drop the trailing paragraph that only carried the old _GoBack
#    bookmark (now relocated above) and append two empty paragraphs.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Delete() | Out-Null

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter() | Out-Null
$end2 = $d.Content
$end2.Collapse(0)
$end2.InsertParagraphAfter() | Out-Null
